$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "UserInput"
$ws.Range("C1").Value = "Expected_Result"
$ws.Range("C1").WrapText = $false

$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

[void]$ws.Range("H8").Select()
